$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute(" user type, and epic)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Found:" $found
$r.InsertBefore(" Importance, Estimate, Acceptance,")
